$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('G2').Value = 'Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Administrator'
$ws.Range('G3').Value = 'Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Majorelle Magdy, Dr. Eman Tantawi, Administrator'
$ws.Range('G4').Value = 'Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Majorelle Magdy, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad'
$ws.Range('G5').Value = 'Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Asmaa Reda, Dr. Veronia Rafat'
$ws.Range('G6').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Manar Montaser, Dr. Majorelle Magdy, Dr. Mohammad El-Tanany, Dr. Alshimaa Atef'
$ws.Range('G7').Value = 'Dr. Amera Ahmad Saad, Dr. Abeer Ragab, Dr. Nada Mohammad, Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Lamiaa Ossama, Dr. Kerelos Zareef'
$ws.Range('G8').Value = 'Dr. Nada Mohammad, Dr. Abeer Ragab'
$ws.Range('G9').Value = 'Dr. Safa Hany, Dr. Shimaa Ashraf'
$ws.Range('G11').Value = 'Dr. Safa Hany, Dr. Amal Awwad, Dr. Aya Saeed'
$ws.Range('G12').Value = 'Dr. Marina Youhanna, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Dina Adel, Dr. Amira Ibrahim, Dr. Eman M. Abo-Sakaya'
$ws.Range('G15').Value = 'Dr. Rania Ahmad Youssef, Dr. Mohammad Safwat'
$ws.Range('G17').Value = 'Dr. Mohammad Safwat, Dr. Esraa Samy'
$ws.Range('G24').Value = 'Dr. Sarah Mahdy, Dr. Youstina Gamil'
$ws.Range('G27').Value = 'Dr. Nourham Mostafa, Dr. Hana Amr'
$ws.Range('G30').Value = 'Dr. Yassmen Ahmad, Dr. Aya Hanafy, Dr. Shorok Mohammad, Dr. Wafaa Ebida'
